$wb = $excel.ActiveWorkbook

# --- Sheet "VEDA_Sets-Comm": rename header + add NRG set refs + split *BIO* out ---
$wsComm = $wb.Worksheets.Item("VEDA_Sets-Comm")

$wsComm.Range("F2").Value = "c_Pos_AndOr"
$wsComm.Range("G2").Value = "c_Neg_AndOr"

$wsComm.Range("A4").Value = "NRG"
$wsComm.Range("F4").Value = "And"
$wsComm.Range("G4").Value = "And"

$wsComm.Range("A7").Value = "NRG"
$wsComm.Range("F7").Value = "And"
$wsComm.Range("G7").Value = "And"

$wsComm.Range("A8").Value = "NRG"
$wsComm.Range("F8").Value = "And"
$wsComm.Range("G8").Value = "And"

$wsComm.Range("A9").Value = "NRG"
$wsComm.Range("B9").Value = "*BIO*"
$wsComm.Range("F9").Value = "Or"
$wsComm.Range("G9").Value = "And"

$wsComm.Range("B8").Value = "*SOL*,*WIN*,*RNW*,*HYD*,-*SOLID*,-ELCRNW,-RNW"

# --- Switch the active sheet / selection to match the saved view state ---
$wsComm.Select()
$wsComm.Range("G3").Select()
